# Add the new "Player Info" worksheet as the first sheet in the workbook,
# shifting "ODI Batting" and "ODI Bowling" one position to the right.
$wb = $excel.ActiveWorkbook

$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($wb.Worksheets.Item(1))

# Header row (bold, centered, thin border - matches the style already used
# for header rows on the other sheets in this workbook).
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row. ID is stored as text (matches how the other sheets store
# numeric-looking identifiers as inline strings).
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5939"
$playerInfo.Range("B2").Value = "Keon Jovani Harding"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# Rename the header + replace the full URL value with just the match code
# on both the "ODI Batting" and "ODI Bowling" sheets. The match code is
# kept as text (matches how the other numeric-looking values on these
# sheets are stored).
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4447"

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4447"
